$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the node version value that changed from v6.3 to v6.9.1 (row 6, column B before the
#    new column is inserted).
$ws.Range("B6").Value = "v6.9.1"

# 2. Insert a new column before column B ("mb version" column).
#    This shifts existing columns B:G to C:H along with their formatting.
$ws.Columns.Item(2).Insert()

# 3. Give the new column the same width as column A.
$ws.Columns.Item(2).ColumnWidth = 19.3333333

# 4. Populate the new column.
$ws.Range("B1").Value = "mb version"
$ws.Range("B2").Value = "EOL v1.7"
$ws.Range("B3").Value = "EOL v1.7"

# 5. Highlight the two EOL rows (v0.10 / v0.12) with an orange fill.
$ws.Range("A2:H3").Interior.Color = 49407

# 6. Update the active selection to match the authored state.
$ws.Range("B4").Select()
